$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) cells keep their original text formatting (no numeric auto-conversion)
$priceRows = @(2,3,5,6,7,9,10,13,14,15,16,17,18,20,21,22,23,24,27,30,31,35,37,38,39,41,42,45,46,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "61.050.05"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.399.61"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "571.28"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "142.22"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "3.398.76"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "7.53"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "3.976.70"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "28.51"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "3.395.96"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "61.103.34"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "14.00"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "9.09"
$ws.Range("E21").Value = "  -4.35%  "
$ws.Range("D22").Value = "386.15"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "74.20"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "3.531.05"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "8.01"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "23.73"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "166.09"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "3.428.13"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").Value = "5.02"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").Value = "28.40"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").Value = "0.0780"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "42.17"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "4.44"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "1.13"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "2.488.23"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("D50").Value = "23.40"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "6.85"
$ws.Range("E51").Value = "  -0.85%  "
